# Updates the "想去人数" (want-to-go count) figures across the workbook's
# four sheets to the newly scraped totals.

$wb = $excel.ActiveWorkbook

# 展览 (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 11331
$ws1.Range("F7").Value  = 306
$ws1.Range("F10").Value = 943
$ws1.Range("F11").Value = 2236
$ws1.Range("F13").Value = 1104
$ws1.Range("F15").Value = 575
$ws1.Range("F16").Value = 853
$ws1.Range("F17").Value = 998
$ws1.Range("F20").Value = 681
$ws1.Range("F21").Value = 709
$ws1.Range("F28").Value = 193
$ws1.Range("F31").Value = 629
$ws1.Range("F32").Value = 2449
$ws1.Range("F35").Value = 121
$ws1.Range("F38").Value = 1500
$ws1.Range("F39").Value = 426
$ws1.Range("F41").Value = 63
$ws1.Range("F43").Value = 54

# 演出 (Performance)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F5").Value  = 77
$ws2.Range("F12").Value = 153
$ws2.Range("F16").Value = 86

# 本地生活 (Local Life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 2219
$ws3.Range("F4").Value = 628

# 全部类型 (All Types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 2219
$ws4.Range("F5").Value  = 628
$ws4.Range("F11").Value = 77
$ws4.Range("F12").Value = 943
$ws4.Range("F13").Value = 2236
$ws4.Range("F15").Value = 1104
$ws4.Range("F17").Value = 575
$ws4.Range("F18").Value = 853
$ws4.Range("F19").Value = 998
$ws4.Range("F23").Value = 681
$ws4.Range("F24").Value = 709
$ws4.Range("F31").Value = 193
$ws4.Range("F34").Value = 2449
$ws4.Range("F35").Value = 153
$ws4.Range("F38").Value = 1500
$ws4.Range("F39").Value = 426
$ws4.Range("F43").Value = 54
